$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A32").Value = "Demo inplannen"
$ws.Range("B32").Value = "klantenservice@testbedrijf123.nl"
$ws.Range("C32").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Range("D32").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E32").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Range("F32").Value = "2025-08-14 21:54:22"
$ws.Range("G32").Value = "Nee"
$ws.Range("H32").Value = "Ja"
$ws.Range("I32").Value = "Nee"
$ws.Range("J32").Value = "Nee"

$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 24

$ws.Range("D2:D31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D32"))
$ws.Range("G2:G31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G32"))
$ws.Range("H2:H31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H32"))
$ws.Range("I2:I31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I32"))
$ws.Range("J2:J31").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J32"))
